$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "61.588.35"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  +0.75%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.391.60"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  -0.18%  "
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  +0.00%  "
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "576.28"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  +0.48%  "
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "141.71"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  -0.43%  "
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.473"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  -0.56%  "
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "7.67"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  +0.47%  "
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  -1.16%  "
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.385"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  -2.33%  "
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "3.972.43"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  -0.12%  "
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  -0.10%  "
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "28.25"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  +1.46%  "
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.420.32"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  +0.66%  "
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  -0.55%  "
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "61.659.23"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  +0.80%  "
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  -0.02%  "
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "13.62"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  -1.18%  "
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "9.02"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  +1.06%  "
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "390.30"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  +1.86%  "
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "74.82"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  +0.09%  "
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  -1.33%  "
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  -0.18%  "
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  -3.92%  "
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.196"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +8.48%  "
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  -0.01%  "
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.99"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  -0.13%  "
$c.Style = "Normal"

$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  -0.70%  "
$c.Style = "Normal"

$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  +1.94%  "
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  -0.02%  "
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "23.31"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  -0.66%  "
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "6.92"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  -1.10%  "
$c.Style = "Normal"

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "169.44"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  +1.71%  "
$c.Style = "Normal"

$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  +0.48%  "
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.425.85"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  -0.12%  "
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  -0.89%  "
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -0.85%  "
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "25.91"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  -4.74%  "
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "4.42"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  +0.07%  "
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  -1.12%  "
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  +1.45%  "
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.449.72"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  -1.18%  "
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "22.69"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  -1.07%  "
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  -2.13%  "
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +0.04%  "
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  -1.23%  "
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  -6.02%  "
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  -1.31%  "
$c.Style = "Normal"

